# Adds two new worksheets (cont_1223_bf_weighted, mmf_ch_weighted) after the
# existing excl_bf_weighted sheet, each populated with a frequency table that
# matches the layout/styling of the first sheet (header row in B1:C1, row
# labels in column A styled the same as the header, numeric Weighted_Count
# data in columns B/C).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Sheet 2: cont_1223_bf_weighted
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "cont_1223_bf_weighted"

# Match the outline/page-setup conventions used by the first sheet.
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws2.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws2.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws2.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws2.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws2.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# Copy the header/label formatting from the first sheet so the new sheet
# reuses the exact same cell style (bold, centered, thin border).
$ws1.Range("B1:C1").Copy()
$ws2.Range("B1:C1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("B1").Value = "cont_1223_bf"
$ws2.Range("C1").Value = "Weighted_Count"

$csv2 = @"
Total|46.1|758.9
CH Sex: Female|45.4|368.3
CH Sex: Male|46.6|390.7
Rural|48.3|526.7
Urban|41|232.2
Central Highlands|61.7|54.1
Mekong River Delta|34.1|113.6
North Central and Central Coast|53|160.5
Northern Midlands and Mountain|57.3|153.6
Red River Delta|44.9|156.3
Southeast|28.4|120.9
Mother Edu: Higher|32.6|97.8
Mother Edu: None/ECE|54.2|33.4
Mother Edu: Primary|45.4|113.4
Mother Edu: Secondary|48.2|514.2
Middle|47|144.3
Poor|49.8|130.7
Poorest|53.7|165.4
Rich|46.5|144.7
Richest|34.8|173.7
Khmer|64.3|10.7
Kinh and Hoa|43|645.6
Mong|62.3|11.6
Other/Missing|76.3|23.5
Tay, Thai, Muong, Nung|58.8|67.5
Elderly HoH: NO|46.1|613.7
Elderly HoH: YES|45.8|145.2
Female|42|185.4
Male|47.4|573.5
"@

$lines2 = $csv2 -split "`n"
$n2 = $lines2.Count
$data2 = New-Object 'object[,]' $n2,3
for ($i = 0; $i -lt $n2; $i++) {
  $parts = $lines2[$i].Split("|")
  $data2[$i,0] = $parts[0]
  $data2[$i,1] = [double]$parts[1]
  $data2[$i,2] = [double]$parts[2]
}
$ws2.Range("A2").Resize($n2, 3).Value = $data2

# ---------------------------------------------------------------------
# Sheet 3: mmf_ch_weighted
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "mmf_ch_weighted"

$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws3.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws3.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws3.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws3.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws3.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

$ws1.Range("B1:C1").Copy()
$ws3.Range("B1:C1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws3.Range("A2:A30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws3.Range("B1").Value = "mmf_ch"
$ws3.Range("C1").Value = "Weighted_Count"

$csv3 = @"
Total|54.6|1099.9
CH Sex: Female|51.6|561
CH Sex: Male|57.6|538.9
Rural|54.4|775
Urban|54.9|324.9
Central Highlands|34.3|75.59999999999999
Mekong River Delta|57.7|183.6
North Central and Central Coast|46.6|224.1
Northern Midlands and Mountain|48.2|218.8
Red River Delta|67.7|221.4
Southeast|61.8|176.5
Mother Edu: Higher|60.2|130.3
Mother Edu: None/ECE|42.1|53.2
Mother Edu: Primary|52.4|167
Mother Edu: Secondary|55|749.3
Middle|48.6|199.2
Poor|52.5|210.4
Poorest|50.4|234.2
Rich|56.9|212.3
Richest|63.2|244
Khmer|47.5|16.7
Kinh and Hoa|56|929.1
Mong|21.2|21.6
Other/Missing|29.2|38.2
Tay, Thai, Muong, Nung|60.2|94.2
Elderly HoH: NO|53.7|875.8
Elderly HoH: YES|58.2|224.1
Female|64.90000000000001|250
Male|51.5|849.9
"@

$lines3 = $csv3 -split "`n"
$n3 = $lines3.Count
$data3 = New-Object 'object[,]' $n3,3
for ($i = 0; $i -lt $n3; $i++) {
  $parts = $lines3[$i].Split("|")
  $data3[$i,0] = $parts[0]
  $data3[$i,1] = [double]$parts[1]
  $data3[$i,2] = [double]$parts[2]
}
$ws3.Range("A2").Resize($n3, 3).Value = $data3

# Restore the original active sheet/selection so the workbook-level view
# state (activeTab) is unchanged by this edit.
$ws1.Activate()
$ws1.Range("A1").Select()
